$wb = $excel.ActiveWorkbook

# ----- Sheet 1 -----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:29:57"
$ws.Cells.Item(3,1).Value = "Total filas: 163"
$arr = New-Object 'object[,]' 163,5
$arr[0,0] = "03:45:25"
$arr[0,1] = "03:47"
$arr[0,2] = "14_ABASTO"
$arr[0,3] = 2
$arr[0,4] = "LP1912"
$arr[1,0] = "03:45:25"
$arr[1,1] = "04:01"
$arr[1,2] = "81_EL PELIGRO"
$arr[1,3] = 16
$arr[1,4] = "LP1912"
$arr[2,0] = "03:45:25"
$arr[2,1] = "04:46"
$arr[2,2] = "215A_EL PATO"
$arr[2,3] = 61
$arr[2,4] = "LP1912"
$arr[3,0] = "03:45:25"
$arr[3,1] = "04:53"
$arr[3,2] = "11_ETCHEVERRY"
$arr[3,3] = 68
$arr[3,4] = "LP1912"
$arr[4,0] = "04:56:49"
$arr[4,1] = "05:13"
$arr[4,2] = "14_ABASTO"
$arr[4,3] = 17
$arr[4,4] = "LP1912"
$arr[5,0] = "03:45:25"
$arr[5,1] = "05:16"
$arr[5,2] = "17_ROMERO"
$arr[5,3] = 91
$arr[5,4] = "LP1912"
$arr[6,0] = "04:45:05"
$arr[6,1] = "05:16"
$arr[6,2] = "14_ABASTO"
$arr[6,3] = 31
$arr[6,4] = "LP1912"
$arr[7,0] = "03:45:25"
$arr[7,1] = "05:22"
$arr[7,2] = "23_HERNANDEZ"
$arr[7,3] = 97
$arr[7,4] = "LP1912"
$arr[8,0] = "05:26:08"
$arr[8,1] = "05:28"
$arr[8,2] = "14_ABASTO"
$arr[8,3] = 2
$arr[8,4] = "LP1912"
$arr[9,0] = "04:18:02"
$arr[9,1] = "05:34"
$arr[9,2] = "14_ABASTO"
$arr[9,3] = 76
$arr[9,4] = "LP1912"
$arr[10,0] = "03:45:25"
$arr[10,1] = "05:34"
$arr[10,2] = "215B_EL PATO"
$arr[10,3] = 109
$arr[10,4] = "LP1912"
$arr[11,0] = "04:18:02"
$arr[11,1] = "05:35"
$arr[11,2] = "215B_EL PATO"
$arr[11,3] = 77
$arr[11,4] = "LP1912"
$arr[12,0] = "03:45:25"
$arr[12,1] = "05:37"
$arr[12,2] = "14_ABASTO"
$arr[12,3] = 112
$arr[12,4] = "LP1912"
$arr[13,0] = "04:18:02"
$arr[13,1] = "05:46"
$arr[13,2] = "15_ABASTO"
$arr[13,3] = 88
$arr[13,4] = "LP1912"
$arr[14,0] = "04:45:05"
$arr[14,1] = "06:04"
$arr[14,2] = "16_SANTA ANA"
$arr[14,3] = 79
$arr[14,4] = "LP1912"
$arr[15,0] = "04:18:02"
$arr[15,1] = "06:05"
$arr[15,2] = "16_SANTA ANA"
$arr[15,3] = 107
$arr[15,4] = "LP1912"
$arr[16,0] = "04:56:49"
$arr[16,1] = "06:11"
$arr[16,2] = "215A_EL PATO"
$arr[16,3] = 75
$arr[16,4] = "LP1912"
$arr[17,0] = "04:18:02"
$arr[17,1] = "06:12"
$arr[17,2] = "215A_EL PATO"
$arr[17,3] = 114
$arr[17,4] = "LP1912"
$arr[18,0] = "04:18:02"
$arr[18,1] = "06:14"
$arr[18,2] = "225_HARAS DEL SUR"
$arr[18,3] = 116
$arr[18,4] = "LP1912"
$arr[19,0] = "04:45:05"
$arr[19,1] = "06:21"
$arr[19,2] = "26_HERNANDEZ"
$arr[19,3] = 96
$arr[19,4] = "LP1912"
$arr[20,0] = "06:25:43"
$arr[20,1] = "06:26"
$arr[20,2] = "86_EST CHICA-ESC AGRARIA"
$arr[20,3] = 1
$arr[20,4] = "LP1912"
$arr[21,0] = "04:45:05"
$arr[21,1] = "06:27"
$arr[21,2] = "23_HERNANDEZ"
$arr[21,3] = 102
$arr[21,4] = "LP1912"
$arr[22,0] = "06:25:43"
$arr[22,1] = "06:28"
$arr[22,2] = "23_HERNANDEZ"
$arr[22,3] = 3
$arr[22,4] = "LP1912"
$arr[23,0] = "04:56:49"
$arr[23,1] = "06:29"
$arr[23,2] = "86_EST CHICA-ESC AGRARIA"
$arr[23,3] = 93
$arr[23,4] = "LP1912"
$arr[24,0] = "04:45:05"
$arr[24,1] = "06:30"
$arr[24,2] = "86_EST CHICA-ESC AGRARIA"
$arr[24,3] = 105
$arr[24,4] = "LP1912"
$arr[25,0] = "04:45:05"
$arr[25,1] = "06:31"
$arr[25,2] = "16_SANTA ANA"
$arr[25,3] = 106
$arr[25,4] = "LP1912"
$arr[26,0] = "04:45:05"
$arr[26,1] = "06:44"
$arr[26,2] = "225_C ROCA-H SUR"
$arr[26,3] = 119
$arr[26,4] = "LP1912"
$arr[27,0] = "05:55:25"
$arr[27,1] = "06:44"
$arr[27,2] = "26_HERNANDEZ"
$arr[27,3] = 49
$arr[27,4] = "LP1912"
$arr[28,0] = "04:56:49"
$arr[28,1] = "06:46"
$arr[28,2] = "215C_EL PATO"
$arr[28,3] = 110
$arr[28,4] = "LP1912"
$arr[29,0] = "05:26:08"
$arr[29,1] = "06:47"
$arr[29,2] = "215C_EL PATO"
$arr[29,3] = 81
$arr[29,4] = "LP1912"
$arr[30,0] = "05:55:25"
$arr[30,1] = "06:59"
$arr[30,2] = "14_ABASTO"
$arr[30,3] = 64
$arr[30,4] = "LP1912"
$arr[31,0] = "05:26:08"
$arr[31,1] = "07:00"
$arr[31,2] = "14_ABASTO"
$arr[31,3] = 94
$arr[31,4] = "LP1912"
$arr[32,0] = "06:25:43"
$arr[32,1] = "07:01"
$arr[32,2] = "16_SANTA ANA"
$arr[32,3] = 36
$arr[32,4] = "LP1912"
$arr[33,0] = "05:55:25"
$arr[33,1] = "07:04"
$arr[33,2] = "23_HERNANDEZ"
$arr[33,3] = 69
$arr[33,4] = "LP1912"
$arr[34,0] = "05:26:08"
$arr[34,1] = "07:05"
$arr[34,2] = "15_ABASTO"
$arr[34,3] = 99
$arr[34,4] = "LP1912"
$arr[35,0] = "05:26:08"
$arr[35,1] = "07:05"
$arr[35,2] = "23_HERNANDEZ"
$arr[35,3] = 99
$arr[35,4] = "LP1912"
$arr[36,0] = "05:26:08"
$arr[36,1] = "07:06"
$arr[36,2] = "10_OLMOS"
$arr[36,3] = 100
$arr[36,4] = "LP1912"
$arr[37,0] = "05:26:08"
$arr[37,1] = "07:07"
$arr[37,2] = "225_GOMEZ"
$arr[37,3] = 101
$arr[37,4] = "LP1912"
$arr[38,0] = "05:26:08"
$arr[38,1] = "07:11"
$arr[38,2] = "215A_EL PATO"
$arr[38,3] = 105
$arr[38,4] = "LP1912"
$arr[39,0] = "06:55:02"
$arr[39,1] = "07:12"
$arr[39,2] = "215A_EL PATO"
$arr[39,3] = 17
$arr[39,4] = "LP1912"
$arr[40,0] = "06:25:43"
$arr[40,1] = "07:14"
$arr[40,2] = "26_HERNANDEZ"
$arr[40,3] = 49
$arr[40,4] = "LP1912"
$arr[41,0] = "05:55:25"
$arr[41,1] = "07:15"
$arr[41,2] = "11_ETCHEVERRY"
$arr[41,3] = 80
$arr[41,4] = "LP1912"
$arr[42,0] = "05:26:08"
$arr[42,1] = "07:16"
$arr[42,2] = "11_ETCHEVERRY"
$arr[42,3] = 110
$arr[42,4] = "LP1912"
$arr[43,0] = "06:55:02"
$arr[43,1] = "07:17"
$arr[43,2] = "16_SANTA ANA"
$arr[43,3] = 22
$arr[43,4] = "LP1912"
$arr[44,0] = "05:26:08"
$arr[44,1] = "07:21"
$arr[44,2] = "26_HERNANDEZ"
$arr[44,3] = 115
$arr[44,4] = "LP1912"
$arr[45,0] = "05:26:08"
$arr[45,1] = "07:23"
$arr[45,2] = "10_OLMOS"
$arr[45,3] = 117
$arr[45,4] = "LP1912"
$arr[46,0] = "05:55:25"
$arr[46,1] = "07:30"
$arr[46,2] = "10_OLMOS"
$arr[46,3] = 95
$arr[46,4] = "LP1912"
$arr[47,0] = "05:55:25"
$arr[47,1] = "07:31"
$arr[47,2] = "16_SANTA ANA"
$arr[47,3] = 96
$arr[47,4] = "LP1912"
$arr[48,0] = "05:55:25"
$arr[48,1] = "07:31"
$arr[48,2] = "11_ETCHEVERRY"
$arr[48,3] = 96
$arr[48,4] = "LP1912"
$arr[49,0] = "05:55:25"
$arr[49,1] = "07:32"
$arr[49,2] = "84_COLONIA URQUIZA-ESC 49"
$arr[49,3] = 97
$arr[49,4] = "LP1912"
$arr[50,0] = "06:55:02"
$arr[50,1] = "07:32"
$arr[50,2] = "16_SANTA ANA"
$arr[50,3] = 37
$arr[50,4] = "LP1912"
$arr[51,0] = "06:55:02"
$arr[51,1] = "07:32"
$arr[51,2] = "11_ETCHEVERRY"
$arr[51,3] = 37
$arr[51,4] = "LP1912"
$arr[52,0] = "07:19:29"
$arr[52,1] = "07:35"
$arr[52,2] = "23_HERNANDEZ"
$arr[52,3] = 16
$arr[52,4] = "LP1912"
$arr[53,0] = "05:55:25"
$arr[53,1] = "07:36"
$arr[53,2] = "27_EL RETIRO"
$arr[53,3] = 101
$arr[53,4] = "LP1912"
$arr[54,0] = "06:55:02"
$arr[54,1] = "07:37"
$arr[54,2] = "27_EL RETIRO"
$arr[54,3] = 42
$arr[54,4] = "LP1912"
$arr[55,0] = "05:55:25"
$arr[55,1] = "07:39"
$arr[55,2] = "10_OLMOS"
$arr[55,3] = 104
$arr[55,4] = "LP1912"
$arr[56,0] = "07:19:29"
$arr[56,1] = "07:46"
$arr[56,2] = "16_SANTA ANA"
$arr[56,3] = 27
$arr[56,4] = "LP1912"
$arr[57,0] = "05:55:25"
$arr[57,1] = "07:47"
$arr[57,2] = "14_ABASTO"
$arr[57,3] = 112
$arr[57,4] = "LP1912"
$arr[58,0] = "06:55:02"
$arr[58,1] = "07:48"
$arr[58,2] = "14_ABASTO"
$arr[58,3] = 53
$arr[58,4] = "LP1912"
$arr[59,0] = "07:50:16"
$arr[59,1] = "07:50"
$arr[59,2] = "10_OLMOS"
$arr[59,3] = 0
$arr[59,4] = "LP1912"
$arr[60,0] = "05:55:25"
$arr[60,1] = "07:51"
$arr[60,2] = "215D_EL PATO"
$arr[60,3] = 116
$arr[60,4] = "LP1912"
$arr[61,0] = "06:55:02"
$arr[61,1] = "07:52"
$arr[61,2] = "215D_EL PATO"
$arr[61,3] = 57
$arr[61,4] = "LP1912"
$arr[62,0] = "07:19:29"
$arr[62,1] = "07:59"
$arr[62,2] = "23_HERNANDEZ"
$arr[62,3] = 40
$arr[62,4] = "LP1912"
$arr[63,0] = "06:25:43"
$arr[63,1] = "08:01"
$arr[63,2] = "23_HERNANDEZ"
$arr[63,3] = 96
$arr[63,4] = "LP1912"
$arr[64,0] = "07:19:29"
$arr[64,1] = "08:03"
$arr[64,2] = "11_ETCHEVERRY"
$arr[64,3] = 44
$arr[64,4] = "LP1912"
$arr[65,0] = "06:55:02"
$arr[65,1] = "08:03"
$arr[65,2] = "23_HERNANDEZ"
$arr[65,3] = 68
$arr[65,4] = "LP1912"
$arr[66,0] = "08:02:22"
$arr[66,1] = "08:05"
$arr[66,2] = "23_HERNANDEZ"
$arr[66,3] = 3
$arr[66,4] = "LP1912"
$arr[67,0] = "07:19:29"
$arr[67,1] = "08:10"
$arr[67,2] = "16_SANTA ANA"
$arr[67,3] = 51
$arr[67,4] = "LP1912"
$arr[68,0] = "07:50:16"
$arr[68,1] = "08:11"
$arr[68,2] = "16_SANTA ANA"
$arr[68,3] = 21
$arr[68,4] = "LP1912"
$arr[69,0] = "06:25:43"
$arr[69,1] = "08:12"
$arr[69,2] = "15_ABASTO"
$arr[69,3] = 107
$arr[69,4] = "LP1912"
$arr[70,0] = "07:50:16"
$arr[70,1] = "08:13"
$arr[70,2] = "10_OLMOS"
$arr[70,3] = 23
$arr[70,4] = "LP1912"
$arr[71,0] = "06:55:02"
$arr[71,1] = "08:21"
$arr[71,2] = "26_HERNANDEZ"
$arr[71,3] = 86
$arr[71,4] = "LP1912"
$arr[72,0] = "06:25:43"
$arr[72,1] = "08:22"
$arr[72,2] = "16_P MOR-SANTA ANA"
$arr[72,3] = 117
$arr[72,4] = "LP1912"
$arr[73,0] = "06:55:02"
$arr[73,1] = "08:23"
$arr[73,2] = "16_P MOR-SANTA ANA"
$arr[73,3] = 88
$arr[73,4] = "LP1912"
$arr[74,0] = "06:25:43"
$arr[74,1] = "08:23"
$arr[74,2] = "215B_EL PATO"
$arr[74,3] = 118
$arr[74,4] = "LP1912"
$arr[75,0] = "06:55:02"
$arr[75,1] = "08:27"
$arr[75,2] = "84_COLONIA URQUIZA-ESC 49"
$arr[75,3] = 92
$arr[75,4] = "LP1912"
$arr[76,0] = "07:50:16"
$arr[76,1] = "08:30"
$arr[76,2] = "23_HERNANDEZ"
$arr[76,3] = 40
$arr[76,4] = "LP1912"
$arr[77,0] = "08:02:22"
$arr[77,1] = "08:33"
$arr[77,2] = "10_OLMOS"
$arr[77,3] = 31
$arr[77,4] = "LP1912"
$arr[78,0] = "08:02:22"
$arr[78,1] = "08:34"
$arr[78,2] = "23_HERNANDEZ"
$arr[78,3] = 32
$arr[78,4] = "LP1912"
$arr[79,0] = "08:32:09"
$arr[79,1] = "08:37"
$arr[79,2] = "23_HERNANDEZ"
$arr[79,3] = 5
$arr[79,4] = "LP1912"
$arr[80,0] = "06:55:02"
$arr[80,1] = "08:42"
$arr[80,2] = "81_EL PELIGRO"
$arr[80,3] = 107
$arr[80,4] = "LP1912"
$arr[81,0] = "07:19:29"
$arr[81,1] = "08:43"
$arr[81,2] = "14_ABASTO"
$arr[81,3] = 84
$arr[81,4] = "LP1912"
$arr[82,0] = "07:50:16"
$arr[82,1] = "08:44"
$arr[82,2] = "14_ABASTO"
$arr[82,3] = 54
$arr[82,4] = "LP1912"
$arr[83,0] = "08:32:09"
$arr[83,1] = "08:53"
$arr[83,2] = "10_OLMOS"
$arr[83,3] = 21
$arr[83,4] = "LP1912"
$arr[84,0] = "06:55:02"
$arr[84,1] = "08:54"
$arr[84,2] = "17_ROMERO"
$arr[84,3] = 119
$arr[84,4] = "LP1912"
$arr[85,0] = "07:19:29"
$arr[85,1] = "09:01"
$arr[85,2] = "215A_EL PATO"
$arr[85,3] = 102
$arr[85,4] = "LP1912"
$arr[86,0] = "07:50:16"
$arr[86,1] = "09:02"
$arr[86,2] = "215A_EL PATO"
$arr[86,3] = 72
$arr[86,4] = "LP1912"
$arr[87,0] = "08:02:22"
$arr[87,1] = "09:03"
$arr[87,2] = "11_ETCHEVERRY"
$arr[87,3] = 61
$arr[87,4] = "LP1912"
$arr[88,0] = "08:32:09"
$arr[88,1] = "09:04"
$arr[88,2] = "11_ETCHEVERRY"
$arr[88,3] = 32
$arr[88,4] = "LP1912"
$arr[89,0] = "08:32:09"
$arr[89,1] = "09:05"
$arr[89,2] = "23_HERNANDEZ"
$arr[89,3] = 33
$arr[89,4] = "LP1912"
$arr[90,0] = "07:19:29"
$arr[90,1] = "09:10"
$arr[90,2] = "16_P MOR-SANTA ANA"
$arr[90,3] = 111
$arr[90,4] = "LP1912"
$arr[91,0] = "07:50:16"
$arr[91,1] = "09:11"
$arr[91,2] = "16_P MOR-SANTA ANA"
$arr[91,3] = 81
$arr[91,4] = "LP1912"
$arr[92,0] = "08:32:09"
$arr[92,1] = "09:13"
$arr[92,2] = "10_OLMOS"
$arr[92,3] = 41
$arr[92,4] = "LP1912"
$arr[93,0] = "07:19:29"
$arr[93,1] = "09:16"
$arr[93,2] = "27_EL RETIRO"
$arr[93,3] = 117
$arr[93,4] = "LP1912"
$arr[94,0] = "07:50:16"
$arr[94,1] = "09:17"
$arr[94,2] = "27_EL RETIRO"
$arr[94,3] = 87
$arr[94,4] = "LP1912"
$arr[95,0] = "07:50:16"
$arr[95,1] = "09:21"
$arr[95,2] = "26_HERNANDEZ"
$arr[95,3] = 91
$arr[95,4] = "LP1912"
$arr[96,0] = "08:02:22"
$arr[96,1] = "09:22"
$arr[96,2] = "16_SANTA ANA"
$arr[96,3] = 80
$arr[96,4] = "LP1912"
$arr[97,0] = "08:02:22"
$arr[97,1] = "09:23"
$arr[97,2] = "11_ETCHEVERRY"
$arr[97,3] = 81
$arr[97,4] = "LP1912"
$arr[98,0] = "08:32:09"
$arr[98,1] = "09:23"
$arr[98,2] = "16_SANTA ANA"
$arr[98,3] = 51
$arr[98,4] = "LP1912"
$arr[99,0] = "07:50:16"
$arr[99,1] = "09:23"
$arr[99,2] = "17_ROMERO"
$arr[99,3] = 93
$arr[99,4] = "LP1912"
$arr[100,0] = "07:50:16"
$arr[100,1] = "09:24"
$arr[100,2] = "11_ETCHEVERRY"
$arr[100,3] = 94
$arr[100,4] = "LP1912"
$arr[101,0] = "07:50:16"
$arr[101,1] = "09:28"
$arr[101,2] = "16_SANTA ANA"
$arr[101,3] = 98
$arr[101,4] = "LP1912"
$arr[102,0] = "07:50:16"
$arr[102,1] = "09:32"
$arr[102,2] = "15_ABASTO"
$arr[102,3] = 102
$arr[102,4] = "LP1912"
$arr[103,0] = "07:50:16"
$arr[103,1] = "09:33"
$arr[103,2] = "10_OLMOS"
$arr[103,3] = 103
$arr[103,4] = "LP1912"
$arr[104,0] = "08:56:29"
$arr[104,1] = "09:34"
$arr[104,2] = "23_HERNANDEZ"
$arr[104,3] = 38
$arr[104,4] = "LP1912"
$arr[105,0] = "08:56:29"
$arr[105,1] = "09:34"
$arr[105,2] = "16_SANTA ANA"
$arr[105,3] = 38
$arr[105,4] = "LP1912"
$arr[106,0] = "08:32:09"
$arr[106,1] = "09:35"
$arr[106,2] = "16_SANTA ANA"
$arr[106,3] = 63
$arr[106,4] = "LP1912"
$arr[107,0] = "08:48:08"
$arr[107,1] = "09:35"
$arr[107,2] = "23_HERNANDEZ"
$arr[107,3] = 47
$arr[107,4] = "LP1912"
$arr[108,0] = "09:35:26"
$arr[108,1] = "09:39"
$arr[108,2] = "23_HERNANDEZ"
$arr[108,3] = 4
$arr[108,4] = "LP1912"
$arr[109,0] = "07:50:16"
$arr[109,1] = "09:42"
$arr[109,2] = "215C_EL PATO"
$arr[109,3] = 112
$arr[109,4] = "LP1912"
$arr[110,0] = "08:02:22"
$arr[110,1] = "09:43"
$arr[110,2] = "14_ABASTO"
$arr[110,3] = 101
$arr[110,4] = "LP1912"
$arr[111,0] = "07:50:16"
$arr[111,1] = "09:44"
$arr[111,2] = "14_ABASTO"
$arr[111,3] = 114
$arr[111,4] = "LP1912"
$arr[112,0] = "09:35:26"
$arr[112,1] = "09:46"
$arr[112,2] = "16_SANTA ANA"
$arr[112,3] = 11
$arr[112,4] = "LP1912"
$arr[113,0] = "08:32:09"
$arr[113,1] = "09:52"
$arr[113,2] = "15_ABASTO"
$arr[113,3] = 80
$arr[113,4] = "LP1912"
$arr[114,0] = "08:56:29"
$arr[114,1] = "09:53"
$arr[114,2] = "10_OLMOS"
$arr[114,3] = 57
$arr[114,4] = "LP1912"
$arr[115,0] = "09:35:26"
$arr[115,1] = "09:58"
$arr[115,2] = "16_SANTA ANA"
$arr[115,3] = 23
$arr[115,4] = "LP1912"
$arr[116,0] = "09:35:26"
$arr[116,1] = "10:03"
$arr[116,2] = "11_ETCHEVERRY"
$arr[116,3] = 28
$arr[116,4] = "LP1912"
$arr[117,0] = "08:56:29"
$arr[117,1] = "10:10"
$arr[117,2] = "16_P MOR-SANTA ANA"
$arr[117,3] = 74
$arr[117,4] = "LP1912"
$arr[118,0] = "08:32:09"
$arr[118,1] = "10:11"
$arr[118,2] = "16_P MOR-SANTA ANA"
$arr[118,3] = 99
$arr[118,4] = "LP1912"
$arr[119,0] = "09:35:26"
$arr[119,1] = "10:12"
$arr[119,2] = "15_ABASTO"
$arr[119,3] = 37
$arr[119,4] = "LP1912"
$arr[120,0] = "09:35:26"
$arr[120,1] = "10:13"
$arr[120,2] = "10_OLMOS"
$arr[120,3] = 38
$arr[120,4] = "LP1912"
$arr[121,0] = "08:32:09"
$arr[121,1] = "10:21"
$arr[121,2] = "26_HERNANDEZ"
$arr[121,3] = 109
$arr[121,4] = "LP1912"
$arr[122,0] = "08:32:09"
$arr[122,1] = "10:22"
$arr[122,2] = "17_ROMERO"
$arr[122,3] = 110
$arr[122,4] = "LP1912"
$arr[123,0] = "09:35:26"
$arr[123,1] = "10:23"
$arr[123,2] = "11_ETCHEVERRY"
$arr[123,3] = 48
$arr[123,4] = "LP1912"
$arr[124,0] = "08:56:29"
$arr[124,1] = "10:26"
$arr[124,2] = "215A_EL PATO"
$arr[124,3] = 90
$arr[124,4] = "LP1912"
$arr[125,0] = "08:32:09"
$arr[125,1] = "10:27"
$arr[125,2] = "215A_EL PATO"
$arr[125,3] = 115
$arr[125,4] = "LP1912"
$arr[126,0] = "10:29:57"
$arr[126,1] = "10:29"
$arr[126,2] = "16_SANTA ANA"
$arr[126,3] = 0
$arr[126,4] = "LP1912"
$arr[127,0] = "10:29:57"
$arr[127,1] = "10:31"
$arr[127,2] = "10_OLMOS"
$arr[127,3] = 2
$arr[127,4] = "LP1912"
$arr[128,0] = "09:35:26"
$arr[128,1] = "10:34"
$arr[128,2] = "23_HERNANDEZ"
$arr[128,3] = 59
$arr[128,4] = "LP1912"
$arr[129,0] = "10:29:57"
$arr[129,1] = "10:34"
$arr[129,2] = "16_SANTA ANA"
$arr[129,3] = 5
$arr[129,4] = "LP1912"
$arr[130,0] = "10:29:57"
$arr[130,1] = "10:39"
$arr[130,2] = "23_HERNANDEZ"
$arr[130,3] = 10
$arr[130,4] = "LP1912"
$arr[131,0] = "10:29:57"
$arr[131,1] = "10:41"
$arr[131,2] = "17_ROMERO"
$arr[131,3] = 12
$arr[131,4] = "LP1912"
$arr[132,0] = "08:48:08"
$arr[132,1] = "10:42"
$arr[132,2] = "17_ROMERO"
$arr[132,3] = 114
$arr[132,4] = "LP1912"
$arr[133,0] = "08:56:29"
$arr[133,1] = "10:43"
$arr[133,2] = "14_ABASTO"
$arr[133,3] = 107
$arr[133,4] = "LP1912"
$arr[134,0] = "08:48:08"
$arr[134,1] = "10:44"
$arr[134,2] = "14_ABASTO"
$arr[134,3] = 116
$arr[134,4] = "LP1912"
$arr[135,0] = "10:29:57"
$arr[135,1] = "10:51"
$arr[135,2] = "15_ABASTO"
$arr[135,3] = 22
$arr[135,4] = "LP1912"
$arr[136,0] = "10:29:57"
$arr[136,1] = "10:52"
$arr[136,2] = "10_OLMOS"
$arr[136,3] = 23
$arr[136,4] = "LP1912"
$arr[137,0] = "09:35:26"
$arr[137,1] = "10:54"
$arr[137,2] = "27_EL RETIRO"
$arr[137,3] = 79
$arr[137,4] = "LP1912"
$arr[138,0] = "10:29:57"
$arr[138,1] = "10:56"
$arr[138,2] = "27_EL RETIRO"
$arr[138,3] = 27
$arr[138,4] = "LP1912"
$arr[139,0] = "10:29:57"
$arr[139,1] = "11:01"
$arr[139,2] = "215C_EL PATO"
$arr[139,3] = 32
$arr[139,4] = "LP1912"
$arr[140,0] = "09:35:26"
$arr[140,1] = "11:02"
$arr[140,2] = "215C_EL PATO"
$arr[140,3] = 87
$arr[140,4] = "LP1912"
$arr[141,0] = "10:29:57"
$arr[141,1] = "11:03"
$arr[141,2] = "11_ETCHEVERRY"
$arr[141,3] = 34
$arr[141,4] = "LP1912"
$arr[142,0] = "10:29:57"
$arr[142,1] = "11:04"
$arr[142,2] = "23_HERNANDEZ"
$arr[142,3] = 35
$arr[142,4] = "LP1912"
$arr[143,0] = "09:35:26"
$arr[143,1] = "11:06"
$arr[143,2] = "16_P MOR-167 Y 521"
$arr[143,3] = 91
$arr[143,4] = "LP1912"
$arr[144,0] = "10:29:57"
$arr[144,1] = "11:11"
$arr[144,2] = "15_ABASTO"
$arr[144,3] = 42
$arr[144,4] = "LP1912"
$arr[145,0] = "09:35:26"
$arr[145,1] = "11:19"
$arr[145,2] = "86_EST CHICA-ESC AGRARIA"
$arr[145,3] = 104
$arr[145,4] = "LP1912"
$arr[146,0] = "10:29:57"
$arr[146,1] = "11:20"
$arr[146,2] = "26_HERNANDEZ"
$arr[146,3] = 51
$arr[146,4] = "LP1912"
$arr[147,0] = "09:35:26"
$arr[147,1] = "11:21"
$arr[147,2] = "26_HERNANDEZ"
$arr[147,3] = 106
$arr[147,4] = "LP1912"
$arr[148,0] = "10:29:57"
$arr[148,1] = "11:26"
$arr[148,2] = "225_C ROCA-H SUR"
$arr[148,3] = 57
$arr[148,4] = "LP1912"
$arr[149,0] = "09:35:26"
$arr[149,1] = "11:27"
$arr[149,2] = "225_C ROCA-H SUR"
$arr[149,3] = 112
$arr[149,4] = "LP1912"
$arr[150,0] = "10:29:57"
$arr[150,1] = "11:31"
$arr[150,2] = "81_EL PELIGRO"
$arr[150,3] = 62
$arr[150,4] = "LP1912"
$arr[151,0] = "09:35:26"
$arr[151,1] = "11:32"
$arr[151,2] = "81_EL PELIGRO"
$arr[151,3] = 117
$arr[151,4] = "LP1912"
$arr[152,0] = "10:29:57"
$arr[152,1] = "11:35"
$arr[152,2] = "11_ETCHEVERRY"
$arr[152,3] = 66
$arr[152,4] = "LP1912"
$arr[153,0] = "10:29:57"
$arr[153,1] = "11:40"
$arr[153,2] = "10_OLMOS"
$arr[153,3] = 71
$arr[153,4] = "LP1912"
$arr[154,0] = "10:29:57"
$arr[154,1] = "11:41"
$arr[154,2] = "17_ROMERO"
$arr[154,3] = 72
$arr[154,4] = "LP1912"
$arr[155,0] = "10:29:57"
$arr[155,1] = "11:50"
$arr[155,2] = "215B_EL PATO"
$arr[155,3] = 81
$arr[155,4] = "LP1912"
$arr[156,0] = "10:29:57"
$arr[156,1] = "11:58"
$arr[156,2] = "225_GOMEZ"
$arr[156,3] = 89
$arr[156,4] = "LP1912"
$arr[157,0] = "10:29:57"
$arr[157,1] = "12:01"
$arr[157,2] = "84_COLONIA URQUIZA-ESC 49"
$arr[157,3] = 92
$arr[157,4] = "LP1912"
$arr[158,0] = "10:29:57"
$arr[158,1] = "12:06"
$arr[158,2] = "16_P MOR-SANTA ANA"
$arr[158,3] = 97
$arr[158,4] = "LP1912"
$arr[159,0] = "10:29:57"
$arr[159,1] = "12:13"
$arr[159,2] = "17_ROMERO"
$arr[159,3] = 104
$arr[159,4] = "LP1912"
$arr[160,0] = "10:29:57"
$arr[160,1] = "12:15"
$arr[160,2] = "14_ABASTO"
$arr[160,3] = 106
$arr[160,4] = "LP1912"
$arr[161,0] = "10:29:57"
$arr[161,1] = "12:20"
$arr[161,2] = "215A_EL PATO"
$arr[161,3] = 111
$arr[161,4] = "LP1912"
$arr[162,0] = "10:29:57"
$arr[162,1] = "12:20"
$arr[162,2] = "26_HERNANDEZ"
$arr[162,3] = 111
$arr[162,4] = "LP1912"
$ws.Range("A6:E168").Value2 = $arr

# ----- Sheet 2 -----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:29:57"
$ws.Cells.Item(3,1).Value = "Total filas: 21"
$arr = New-Object 'object[,]' 21,5
$arr[0,0] = "03:45:25"
$arr[0,1] = "04:46"
$arr[0,2] = "215A_EL PATO"
$arr[0,3] = 61
$arr[0,4] = "LP1912"
$arr[1,0] = "03:45:25"
$arr[1,1] = "05:34"
$arr[1,2] = "215B_EL PATO"
$arr[1,3] = 109
$arr[1,4] = "LP1912"
$arr[2,0] = "04:18:02"
$arr[2,1] = "05:35"
$arr[2,2] = "215B_EL PATO"
$arr[2,3] = 77
$arr[2,4] = "LP1912"
$arr[3,0] = "04:56:49"
$arr[3,1] = "06:11"
$arr[3,2] = "215A_EL PATO"
$arr[3,3] = 75
$arr[3,4] = "LP1912"
$arr[4,0] = "04:18:02"
$arr[4,1] = "06:12"
$arr[4,2] = "215A_EL PATO"
$arr[4,3] = 114
$arr[4,4] = "LP1912"
$arr[5,0] = "04:56:49"
$arr[5,1] = "06:46"
$arr[5,2] = "215C_EL PATO"
$arr[5,3] = 110
$arr[5,4] = "LP1912"
$arr[6,0] = "05:26:08"
$arr[6,1] = "06:47"
$arr[6,2] = "215C_EL PATO"
$arr[6,3] = 81
$arr[6,4] = "LP1912"
$arr[7,0] = "05:26:08"
$arr[7,1] = "07:11"
$arr[7,2] = "215A_EL PATO"
$arr[7,3] = 105
$arr[7,4] = "LP1912"
$arr[8,0] = "06:55:02"
$arr[8,1] = "07:12"
$arr[8,2] = "215A_EL PATO"
$arr[8,3] = 17
$arr[8,4] = "LP1912"
$arr[9,0] = "05:55:25"
$arr[9,1] = "07:51"
$arr[9,2] = "215D_EL PATO"
$arr[9,3] = 116
$arr[9,4] = "LP1912"
$arr[10,0] = "06:55:02"
$arr[10,1] = "07:52"
$arr[10,2] = "215D_EL PATO"
$arr[10,3] = 57
$arr[10,4] = "LP1912"
$arr[11,0] = "06:25:43"
$arr[11,1] = "08:23"
$arr[11,2] = "215B_EL PATO"
$arr[11,3] = 118
$arr[11,4] = "LP1912"
$arr[12,0] = "07:19:29"
$arr[12,1] = "09:01"
$arr[12,2] = "215A_EL PATO"
$arr[12,3] = 102
$arr[12,4] = "LP1912"
$arr[13,0] = "07:50:16"
$arr[13,1] = "09:02"
$arr[13,2] = "215A_EL PATO"
$arr[13,3] = 72
$arr[13,4] = "LP1912"
$arr[14,0] = "07:50:16"
$arr[14,1] = "09:42"
$arr[14,2] = "215C_EL PATO"
$arr[14,3] = 112
$arr[14,4] = "LP1912"
$arr[15,0] = "08:56:29"
$arr[15,1] = "10:26"
$arr[15,2] = "215A_EL PATO"
$arr[15,3] = 90
$arr[15,4] = "LP1912"
$arr[16,0] = "08:32:09"
$arr[16,1] = "10:27"
$arr[16,2] = "215A_EL PATO"
$arr[16,3] = 115
$arr[16,4] = "LP1912"
$arr[17,0] = "10:29:57"
$arr[17,1] = "11:01"
$arr[17,2] = "215C_EL PATO"
$arr[17,3] = 32
$arr[17,4] = "LP1912"
$arr[18,0] = "09:35:26"
$arr[18,1] = "11:02"
$arr[18,2] = "215C_EL PATO"
$arr[18,3] = 87
$arr[18,4] = "LP1912"
$arr[19,0] = "10:29:57"
$arr[19,1] = "11:50"
$arr[19,2] = "215B_EL PATO"
$arr[19,3] = 81
$arr[19,4] = "LP1912"
$arr[20,0] = "10:29:57"
$arr[20,1] = "12:20"
$arr[20,2] = "215A_EL PATO"
$arr[20,3] = 111
$arr[20,4] = "LP1912"
$ws.Range("A6:E26").Value2 = $arr

# ----- Sheet 3 -----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 10:29:57"
$ws.Cells.Item(3,1).Value = "Total filas: 29"
$arr = New-Object 'object[,]' 29,5
$arr[0,0] = "04:56:49"
$arr[0,1] = "05:43"
$arr[0,2] = "215A_LA PLATA"
$arr[0,3] = 47
$arr[0,4] = "L6173"
$arr[1,0] = "03:45:25"
$arr[1,1] = "05:44"
$arr[1,2] = "215A_LA PLATA"
$arr[1,3] = 119
$arr[1,4] = "L6173"
$arr[2,0] = "04:56:49"
$arr[2,1] = "06:08"
$arr[2,2] = "215A_LA PLATA"
$arr[2,3] = 72
$arr[2,4] = "L6173"
$arr[3,0] = "04:18:02"
$arr[3,1] = "06:09"
$arr[3,2] = "215A_LA PLATA"
$arr[3,3] = 111
$arr[3,4] = "L6173"
$arr[4,0] = "04:56:49"
$arr[4,1] = "06:32"
$arr[4,2] = "215C_LA PLATA"
$arr[4,3] = 96
$arr[4,4] = "L6203"
$arr[5,0] = "04:45:05"
$arr[5,1] = "06:33"
$arr[5,2] = "215C_LA PLATA"
$arr[5,3] = 108
$arr[5,4] = "L6203"
$arr[6,0] = "06:25:43"
$arr[6,1] = "06:59"
$arr[6,2] = "215B_LP-P MOR-1 Y 57"
$arr[6,3] = 34
$arr[6,4] = "L6173"
$arr[7,0] = "05:26:08"
$arr[7,1] = "07:00"
$arr[7,2] = "215B_LP-P MOR-1 Y 57"
$arr[7,3] = 94
$arr[7,4] = "L6173"
$arr[8,0] = "05:55:25"
$arr[8,1] = "07:35"
$arr[8,2] = "215A_LA PLATA"
$arr[8,3] = 100
$arr[8,4] = "L6173"
$arr[9,0] = "06:25:43"
$arr[9,1] = "07:39"
$arr[9,2] = "215A_LA PLATA"
$arr[9,3] = 74
$arr[9,4] = "L6173"
$arr[10,0] = "06:55:02"
$arr[10,1] = "07:42"
$arr[10,2] = "215A_LA PLATA"
$arr[10,3] = 47
$arr[10,4] = "L6173"
$arr[11,0] = "07:19:29"
$arr[11,1] = "07:46"
$arr[11,2] = "215A_LA PLATA"
$arr[11,3] = 27
$arr[11,4] = "L6173"
$arr[12,0] = "07:50:16"
$arr[12,1] = "07:51"
$arr[12,2] = "215A_LA PLATA"
$arr[12,3] = 1
$arr[12,4] = "L6173"
$arr[13,0] = "06:25:43"
$arr[13,1] = "08:06"
$arr[13,2] = "215C_LA PLATA"
$arr[13,3] = 101
$arr[13,4] = "L6203"
$arr[14,0] = "06:55:02"
$arr[14,1] = "08:07"
$arr[14,2] = "215C_LA PLATA"
$arr[14,3] = 72
$arr[14,4] = "L6203"
$arr[15,0] = "07:19:29"
$arr[15,1] = "08:21"
$arr[15,2] = "215C_LA PLATA"
$arr[15,3] = 62
$arr[15,4] = "L6203"
$arr[16,0] = "08:02:22"
$arr[16,1] = "08:25"
$arr[16,2] = "215C_LA PLATA"
$arr[16,3] = 23
$arr[16,4] = "L6203"
$arr[17,0] = "07:50:16"
$arr[17,1] = "08:27"
$arr[17,2] = "215C_LA PLATA"
$arr[17,3] = 37
$arr[17,4] = "L6203"
$arr[18,0] = "07:19:29"
$arr[18,1] = "08:35"
$arr[18,2] = "215A_LA PLATA"
$arr[18,3] = 76
$arr[18,4] = "L6173"
$arr[19,0] = "06:55:02"
$arr[19,1] = "08:36"
$arr[19,2] = "215A_LA PLATA"
$arr[19,3] = 101
$arr[19,4] = "L6173"
$arr[20,0] = "08:02:22"
$arr[20,1] = "08:39"
$arr[20,2] = "215A_LA PLATA"
$arr[20,3] = 37
$arr[20,4] = "L6173"
$arr[21,0] = "08:32:09"
$arr[21,1] = "08:42"
$arr[21,2] = "215A_LA PLATA"
$arr[21,3] = 10
$arr[21,4] = "L6173"
$arr[22,0] = "07:19:29"
$arr[22,1] = "09:09"
$arr[22,2] = "215D_LA PLATA"
$arr[22,3] = 110
$arr[22,4] = "L6203"
$arr[23,0] = "08:32:09"
$arr[23,1] = "10:03"
$arr[23,2] = "215B_LP-P MOR-40 Y 115"
$arr[23,3] = 91
$arr[23,4] = "L6173"
$arr[24,0] = "10:29:57"
$arr[24,1] = "10:53"
$arr[24,2] = "215A_LA PLATA"
$arr[24,3] = 24
$arr[24,4] = "L6173"
$arr[25,0] = "08:56:29"
$arr[25,1] = "10:54"
$arr[25,2] = "215A_LA PLATA"
$arr[25,3] = 118
$arr[25,4] = "L6173"
$arr[26,0] = "10:29:57"
$arr[26,1] = "11:13"
$arr[26,2] = "215C_LA PLATA"
$arr[26,3] = 44
$arr[26,4] = "L6203"
$arr[27,0] = "09:35:26"
$arr[27,1] = "11:14"
$arr[27,2] = "215C_LA PLATA"
$arr[27,3] = 99
$arr[27,4] = "L6203"
$arr[28,0] = "10:29:57"
$arr[28,1] = "12:03"
$arr[28,2] = "215A_LA PLATA"
$arr[28,3] = 94
$arr[28,4] = "L6173"
$ws.Range("A6:E34").Value2 = $arr
